$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.134489417076111
$ws.Range("B1").Value = 1.882919430732727
$ws.Range("D1").Value = 1.906922340393066
$ws.Range("E1").Value = 1.096083402633667
